# Insert a new title paragraph "PARTIE 1 :" (bold, 18pt) before the
# current first paragraph of the document.

$d = $word.ActiveDocument

# Insert the text at the very start of the document (no paragraph mark
# yet), so it does not inherit any character formatting - such as the
# superscript used later in the document - from neighbouring runs.
$start = $d.Range(0, 0)
$start.InsertBefore("PARTIE 1 :")

# Split the text we just typed into its own paragraph by inserting a
# paragraph mark right after it.
$titleText = $d.Range(0, 10)
$titleText.InsertParagraphAfter()

# Apply bold + 18pt (sz 36 half-points) to the whole new paragraph,
# including its paragraph mark, so both the run and the paragraph mark
# pick up the formatting.
$titlePara = $d.Paragraphs(1).Range
$titlePara.Font.Bold = 1
$titlePara.Font.Size = 18
